$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the pre-existing "_GoBack" bookmark. It currently wraps part
#    of the "angular-tour-of-heroes-backend - Copy (2) before signIn"
#    line; the edit relocates the bookmark to the newly inserted
#    "nodemon start" line further up, so the old one must go away.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Find the paragraph that reads
#       "C:\Users\mecostantino\angular>cd nodejs-express-mongodb"
#    (split across several runs/proofErr tags) so we can anchor the
#    insertion point right after it, regardless of its exact index.
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*angular>cd nodejs-express-mongodb*") {
        $targetPara = $d.Paragraphs.Item($i)
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate anchor paragraph"
}

# Collapse a whole-document Range to the character position right after
# the anchor paragraph's text (but before its paragraph mark) so the
# new paragraphs land between it and the following (empty) paragraph.
$insertAt = $targetPara.Range.End
$r = $d.Range($insertAt, $insertAt)

# ------------------------------------------------------------------
# 3) Insert the two new paragraphs ("Then " / the nodemon command line
#    with the relocated _GoBack bookmark). A throw-away trailing
#    <w:p/> is appended to the fragment because InsertXML always fuses
#    the last inserted paragraph with whatever paragraph used to sit
#    at the insertion point -- the extra paragraph absorbs that fusion
#    so the original blank paragraph survives as its own element.
# ------------------------------------------------------------------
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r><w:t xml:space="preserve">Then </w:t></w:r>
</w:p>
<w:p>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>:\Users\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>mecostantino</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>\angular\</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>nodejs</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>-express-</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>mongodb</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>nodemon</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> start</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r.InsertXML($xml)

# ------------------------------------------------------------------
# 4) The insertion above leaves two adjacent empty paragraphs where
#    the document originally had only one (the genuinely-blank one
#    that separated the "cd ..." line from the "//this command..."
#    line). Remove the extra one -- it is the paragraph right after
#    the newly inserted "nodemon start" paragraph.
# ------------------------------------------------------------------
$afterIndex = $targetPara.Index + 3
$dup = $d.Paragraphs.Item($afterIndex)
if ($dup.Range.Text.Trim() -eq "") {
    $dup.Range.Delete()
}

Write-Output "Inserted 'Then' / nodemon-start paragraphs and relocated the _GoBack bookmark."
